$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill in row 8 (C8, D8, E8) with the new timesheet entry
$ws.Range("C8").Value = "3:00pm – 4:45"
$ws.Range("D8").Value = 1.75
$ws.Range("E8").Value = "Registration page reconstruction, code documentation, questions stuff"

# Move the active selection to E9, matching the post-edit cursor position
$ws.Range("E9").Select()
